$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.858.16"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.126.90"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'243.46"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'617.71"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("D7").Value = "'1.12"
$ws.Range("E7").Value = "  -3.59%  "
$ws.Range("D8").Value = "'0.392"
$ws.Range("E8").Value = "  +4.53%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "3.123.45"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "'0.759"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "'0.206"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'35.39"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'5.62"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "91.645.34"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D18").Value = "3.127.59"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'3.79"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'14.93"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "'5.90"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'457.46"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  -6.25%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'5.96"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").Value = "'89.88"
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.48"
$ws.Range("E27").Value = "  +46.37%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'11.75"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.143"
$ws.Range("E30").Value = "  +16.81%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  -6.90%  "
$ws.Range("D34").Value = "'9.40"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "'0.173"
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'26.46"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.51"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'2.00"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("D39").Value = "'491.91"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'3.86"
$ws.Range("E41").Value = "  -8.53%  "
$ws.Range("D42").Value = "'0.439"
$ws.Range("E43").Value = "  -6.56%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "'156.81"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'4.46"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").Value = "'0.0327"
$ws.Range("E51").Value = "  +1.47%  "
